$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '69.327.01'
Set-TextValue $ws.Range("E2") '  -0.25%  '
Set-TextValue $ws.Range("D3") '3.432.31'
Set-TextValue $ws.Range("E3") '  +1.17%  '
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  +0.16%  '
Set-TextValue $ws.Range("D5") '583.14'
Set-TextValue $ws.Range("E5") '  -1.09%  '
Set-TextValue $ws.Range("D6") '178.94'
Set-TextValue $ws.Range("E6") '  -0.31%  '
Set-TextValue $ws.Range("E7") '  +0.12%  '
Set-TextValue $ws.Range("E8") '  -0.65%  '
Set-TextValue $ws.Range("D9") '0.198'
Set-TextValue $ws.Range("E9") '  +5.48%  '
Set-TextValue $ws.Range("D10") '0.587'
Set-TextValue $ws.Range("E10") '  -0.17%  '
Set-TextValue $ws.Range("D11") '48.33'
Set-TextValue $ws.Range("E11") '  -0.61%  '
Set-TextValue $ws.Range("D12") '0.0000282'
Set-TextValue $ws.Range("E12") '  +1.60%  '
Set-TextValue $ws.Range("D13") '681.62'
Set-TextValue $ws.Range("E13") '  -2.91%  '
Set-TextValue $ws.Range("D14") '3.976.32'
Set-TextValue $ws.Range("E14") '  +0.95%  '
Set-TextValue $ws.Range("D15") '8.67'
Set-TextValue $ws.Range("E15") '  +1.47%  '
Set-TextValue $ws.Range("D16") '69.353.55'
Set-TextValue $ws.Range("E16") '  +0.00%  '
Set-TextValue $ws.Range("D17") '3.427.01'
Set-TextValue $ws.Range("E17") '  +1.12%  '
Set-TextValue $ws.Range("E18") '  +0.40%  '
Set-TextValue $ws.Range("D19") '17.84'
Set-TextValue $ws.Range("E19") '  +1.06%  '
Set-TextValue $ws.Range("D20") '11.36'
Set-TextValue $ws.Range("E20") '  +0.20%  '
Set-TextValue $ws.Range("D21") '0.913'
Set-TextValue $ws.Range("E21") '  +0.87%  '
Set-TextValue $ws.Range("D22") '5.39'
Set-TextValue $ws.Range("E22") '  -2.62%  '
Set-TextValue $ws.Range("D23") '16.98'
Set-TextValue $ws.Range("E23") '  -1.08%  '
Set-TextValue $ws.Range("D24") '101.15'
Set-TextValue $ws.Range("E24") '  -0.19%  '
Set-TextValue $ws.Range("E25") '  -1.09%  '
Set-TextValue $ws.Range("E26") '  -1.00%  '
Set-TextValue $ws.Range("D27") '9.69'
Set-TextValue $ws.Range("E27") '  +0.16%  '
Set-TextValue $ws.Range("D28") '33.67'
Set-TextValue $ws.Range("E28") '  +0.50%  '
Set-TextValue $ws.Range("D29") '8.79'
Set-TextValue $ws.Range("E29") '  +1.38%  '
Set-TextValue $ws.Range("D30") '6.90'
Set-TextValue $ws.Range("E30") '  -2.27%  '
Set-TextValue $ws.Range("E31") '  +8.04%  '
Set-TextValue $ws.Range("D32") '562.76'
Set-TextValue $ws.Range("E32") '  +1.20%  '
Set-TextValue $ws.Range("D33") '11.05'
Set-TextValue $ws.Range("E33") '  -1.05%  '
Set-TextValue $ws.Range("E34") '  -1.28%  '
Set-TextValue $ws.Range("D35") '58.10'
Set-TextValue $ws.Range("E35") '  -0.43%  '
Set-TextValue $ws.Range("E36") '  +0.09%  '
Set-TextValue $ws.Range("D37") '3.636.25'
Set-TextValue $ws.Range("E37") '  -2.58%  '
Set-TextValue $ws.Range("D38") '0.141'
Set-TextValue $ws.Range("E38") '  -3.58%  '
Set-TextValue $ws.Range("D39") '35.24'
Set-TextValue $ws.Range("E39") '  +0.76%  '
Set-TextValue $ws.Range("D40") '0.0₃0735'
Set-TextValue $ws.Range("E40") '  +7.25%  '
Set-TextValue $ws.Range("E41") '  +1.62%  '
Set-TextValue $ws.Range("D42") '2.70'
Set-TextValue $ws.Range("E42") '  +1.34%  '
Set-TextValue $ws.Range("D43") '3.35'
Set-TextValue $ws.Range("E43") '  +3.63%  '
Set-TextValue $ws.Range("D44") '0.0425'
Set-TextValue $ws.Range("E44") '  +1.38%  '
Set-TextValue $ws.Range("D45") '0.336'
Set-TextValue $ws.Range("E45") '  -1.23%  '
Set-TextValue $ws.Range("E46") '  -0.01%  '
Set-TextValue $ws.Range("D47") '1.42'
Set-TextValue $ws.Range("E48") '  -0.48%  '
Set-TextValue $ws.Range("D49") '0.999'
Set-TextValue $ws.Range("E49") '  +0.10%  '
Set-TextValue $ws.Range("D50") '131.13'
Set-TextValue $ws.Range("E50") '  -0.88%  '
Set-TextValue $ws.Range("D51") '2.71'
Set-TextValue $ws.Range("E51") '  +2.21%  '
